# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rows 3-6, column F
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 279
$wsExhibition.Range("F4").Value = 935
$wsExhibition.Range("F5").Value = 79
$wsExhibition.Range("F6").Value = 49

# Sheet "全部类型" (All Types) - rows 4-7, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 279
$wsAll.Range("F5").Value = 935
$wsAll.Range("F6").Value = 79
$wsAll.Range("F7").Value = 49
